$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Capture the existing hyperlink cell style (column F uses the built-in Hyperlink style)
# before anything is modified, then remove the stale hyperlink relationships so that
# re-adding them below (in row order) produces a clean rId1..rId11 sequence matching
# the shifted rows, instead of stacking duplicate relationships on top of the old ones.
$hlStyle = $ws.Range("F2").Style
$ws.Hyperlinks.Delete()

$rows = @(
    @('2026-01-06 18:28:30', '大手SIer等のAIソリューション開発・導入を支援してくださるエンジニア・PM募集', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5455098', '375', '🔥AI,Ai ◆開発'),
    @('2026-01-06 18:28:30', '法人向け生成AIサービス(RAG・議事録機能)の設計・開発を支援エンジニア募集(AI/バックエンド)', 'システム開発', '200,000 円 ~ 300,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5445159', '368', '🔥AI,Ai ◆開発'),
    @('2026-01-06 18:28:30', 'B2B向け生成AIサービス(チャット・RAG)の新規開発プロジェクト推進を支援してくださるPM募集', 'システム開発', '200,000 円 ~ 300,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5445154', '368', '🔥AI,Ai ◆開発'),
    @('2026-01-06 18:28:30', 'python等を用いたcsvからデータベース、ポータルの構築と指導', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5466190', '193', '🔥Python'),
    @('2026-01-06 18:28:30', '​【1万〜3万円/BASE経験者】アリエク・ネッシー等のCSVをBASE用に変換・加工するツール作成', 'システム開発', '10,000 円 ~ 20,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5465992', '65', '◆ツール'),
    @('2026-01-06 18:28:30', '初回 【急募】ECサイトの要件定義や基本設計ができる方を募集(1人月、フルリモート可、2025年12月〜)', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5425629', '45', '◇サイト'),
    @('2026-01-06 18:28:30', '【インテリア業界向け】マッチングサイト運用サポートスタッフ募集', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5466189', '33', '◇サイト'),
    @('2026-01-06 18:28:30', 'WEBサーバーの管理、トラブル解決対応できる方を募集します!', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5466047', '33', '◇管理'),
    @('2026-01-06 18:28:30', '現行のシステムに追加要素', 'システム開発', '100,000 円 ~ 200,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5465878', '33', $null),
    @('2026-01-06 18:28:30', '《長期レギュラー》公的機関Web運用の要となる、ディレクター募集', 'システム開発', '100,000 円 ~ 200,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5465685', '18', $null),
    @('2026-01-06 18:28:30', 'ActiveDirectoryの移行(フェーズ1)', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5465836', '13', $null)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row[5]) | Out-Null
    $ws.Cells.Item($r, 6).Style = $hlStyle
    $ws.Cells.Item($r, 7).Value = [double]$row[6]
    if ($row[7] -ne $null) {
        $ws.Cells.Item($r, 8).Value = $row[7]
    } else {
        $ws.Cells.Item($r, 8).ClearContents()
    }
}

Write-Output "done"
